$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 3.270036666666666
$ws.Range("H2").Value = 9.81011
$ws.Range("I2").Value = 0.359406393324744
$ws.Range("J2").Value = 0.3594063933247441
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 8.828998666666665
$ws.Range("N2").Value = 26.486996
$ws.Range("O2").Value = 0.1794455804823882
$ws.Range("P2").Value = 0.1794455804823882
$ws.Range("Q2").Value = 28.87114936995111
$ws.Range("R2").Value = 259.84034432956
$ws.Range("S2").Value = 0.06449388887924021
$ws.Range("T2").Value = 0.06449388887924021

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 3.270036666666666
$ws.Range("H3").Value = 9.81011
$ws.Range("I3").Value = 0.359406393324744
$ws.Range("J3").Value = 0.3594063933247441
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 19.33828433333333
$ws.Range("N3").Value = 58.014853
$ws.Range("O3").Value = 0.3930422677296217
$ws.Range("P3").Value = 0.3930422677296217
$ws.Range("Q3").Value = 63.23689884042555
$ws.Range("R3").Value = 569.13208956383
$ws.Range("S3").Value = 0.1412619038688817
$ws.Range("T3").Value = 0.1412619038688818

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 3.270036666666666
$ws.Range("H4").Value = 9.81011
$ws.Range("I4").Value = 0.359406393324744
$ws.Range("J4").Value = 0.3594063933247441
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 21.03425566666667
$ws.Range("N4").Value = 63.102767
$ws.Range("O4").Value = 0.4275121517879902
$ws.Range("P4").Value = 0.4275121517879902
$ws.Range("Q4").Value = 68.78278728604111
$ws.Range("R4").Value = 619.04508557437
$ws.Range("S4").Value = 0.1536506005766221
$ws.Range("T4").Value = 0.1536506005766221

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 5.828401
$ws.Range("H5").Value = 17.485203
$ws.Range("I5").Value = 0.6405936066752559
$ws.Range("J5").Value = 0.640593606675256
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 8.828998666666665
$ws.Range("N5").Value = 26.486996
$ws.Range("O5").Value = 0.1794455804823882
$ws.Range("P5").Value = 0.1794455804823882
$ws.Range("Q5").Value = 51.45894465779866
$ws.Range("R5").Value = 463.130501920188
$ws.Range("S5").Value = 0.1149516916031479
$ws.Range("T5").Value = 0.1149516916031479

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 5.828401
$ws.Range("H6").Value = 17.485203
$ws.Range("I6").Value = 0.6405936066752559
$ws.Range("J6").Value = 0.640593606675256
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 19.33828433333333
$ws.Range("N6").Value = 58.014853
$ws.Range("O6").Value = 0.3930422677296217
$ws.Range("P6").Value = 0.3930422677296217
$ws.Range("Q6").Value = 112.7112757466844
$ws.Range("R6").Value = 1014.401481720159
$ws.Range("S6").Value = 0.2517803638607399
$ws.Range("T6").Value = 0.2517803638607399

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 5.828401
$ws.Range("H7").Value = 17.485203
$ws.Range("I7").Value = 0.6405936066752559
$ws.Range("J7").Value = 0.640593606675256
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 21.03425566666667
$ws.Range("N7").Value = 63.102767
$ws.Range("O7").Value = 0.4275121517879902
$ws.Range("P7").Value = 0.4275121517879902
$ws.Range("Q7").Value = 122.5960767618557
$ws.Range("R7").Value = 1103.364690856701
$ws.Range("S7").Value = 0.2738615512113681
$ws.Range("T7").Value = 0.2738615512113681
